$d = $word.ActiveDocument

# Locate the paragraph that currently ends the document ("Thanks for reading!")
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Thanks for reading!") {
        $target = $p
    }
}
if ($null -eq $target) {
    throw "Could not locate the 'Thanks for reading!' paragraph"
}

$insertPos = $target.Range.End
$insertRange = $d.Range($insertPos, $insertPos)

$xmlFragment = @"
    <w:p>
      <w:pPr>
        <w:pStyle w:val="BodyText"/>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="BodyText"/>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>Use of AI</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="BodyText"/>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">In addition to the adapter, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>ChatGPT</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> was used for </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">consulting on errors, cleaning some of the code </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">and </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>rephrasing some of</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> the learning diary parts.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> Android Studio AI was tried out, but not used for actual work</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> (isn’t very helpful!)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
"@

$null = $insertRange.InsertXML($xmlFragment)

Write-Output "Inserted AI-usage section after 'Thanks for reading!'"
